# "Generate Report for Handoff"
# Updates the localization-status report: flips the per-language status
# from "In Translation" to "Ready for handoff", refreshes the handoff
# timestamps, and widens the "Status"/language-status columns so the new
# (longer) status text fits.

$wb = $excel.ActiveWorkbook

$newStatus   = "Ready for handoff"
$newDateG    = "2016-09-06 05:14:59"   # Overview "Latest HO Xliff Generate Date" / de-de "Latest Handoff Datetime"
$newDateZh   = "2016-09-06 05:14:47"   # zh-cn "Latest Handoff Datetime"

# Target ColumnWidth (character units) that reproduces the wider OOXML
# column width (17.2159881591797) used for the Status columns.
$newColWidth = 16.38265482584637

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("G2").Value = $newDateG
$ws.Columns.Item(5).ColumnWidth = $newColWidth
$ws.Columns.Item(6).ColumnWidth = $newColWidth

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = $newDateZh
$ws.Columns.Item(3).ColumnWidth = $newColWidth

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = $newDateG
$ws.Columns.Item(3).ColumnWidth = $newColWidth
